$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 & 5 values (ordered to match shared-string allocation order) ---
$ws.Range("C4").Value = "Murrelektronik"
$ws.Range("C5").Value = "Murrelektronik"
$ws.Range("L4").Value = "https://www.automationdirect.com/adc/shopping/catalog/power_products_(electrical)/surge_suppression_devices/universal_surge_suppressors/26183"
$ws.Range("L5").Value = "https://www.automationdirect.com/adc/shopping/catalog/power_products_(electrical)/surge_suppression_devices/universal_surge_suppressors/26051"
$ws.Range("B4").Value = "110-250 VAC 0-60 Hz Varistor Module"
$ws.Range("B5").Value = "12-30 VDC Diode-Zender Suppressor Module"

$ws.Range("A4").Value = "Electronics"
$ws.Range("D4").Value = 26183
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.25

$ws.Range("D5").Value = 26051
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8

# --- Formulas ---
$ws.Range("I3:I5").Formula = "=CEILING.MATH(E3/F3)"
$ws.Range("J4").Formula = "=E4*(G4/F4)"
$ws.Range("K4").Formula = "=I4*G4"
$ws.Range("J5").Formula = "=E5*(G5/F5)"
$ws.Range("K5").Formula = "=I5*G5"
